# Auto-generated edit script applying Famfrit_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8217.9
$ws.Range("I19").Value = 1995.5
$ws.Range("J19").Value = 9773.5
$ws.Range("K19").Value = 1995.5
$ws.Range("L19").Value = 9773.5
$ws.Range("M19").Value = -1820.5
$ws.Range("N19").Value = -10123.5
$ws.Range("H28").Value = 3302.7058
$ws.Range("I28").Value = 504.55554
$ws.Range("J28").Value = 6450.625
$ws.Range("K28").Value = 504.55554
$ws.Range("L28").Value = 6450.625
$ws.Range("M28").Value = -19.55554000000001
$ws.Range("N28").Value = -7420.625
$ws.Range("H54").Value = 7038
$ws.Range("I54").Value = 7038
$ws.Range("K54").Value = 7038
$ws.Range("M54").Value = -6552
$ws.Range("H107").Value = 810.9
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 621.8
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 621.8
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -4461.8
$ws.Range("H135").Value = 7577340
$ws.Range("I135").Value = 1075.9565
$ws.Range("K135").Value = 9683.6085
$ws.Range("M135").Value = -7148.6085
$ws.Range("H138").Value = 4453723
$ws.Range("J138").Value = 4840521
$ws.Range("L138").Value = 14521563
$ws.Range("N138").Value = -14531843

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19242384
$ws.Range("I32").Value = 21747870
$ws.Range("K32").Value = 21747870
$ws.Range("M32").Value = -21747583
$ws.Range("H61").Value = 20412442
$ws.Range("I61").Value = 24394088
$ws.Range("J61").Value = 6499.25
$ws.Range("K61").Value = 24394088
$ws.Range("L61").Value = 6499.25
$ws.Range("M61").Value = -24393876
$ws.Range("N61").Value = -6923.25
$ws.Range("H132").Value = 27786218
$ws.Range("I132").Value = 9713.267
$ws.Range("J132").Value = 166668740
$ws.Range("K132").Value = 29139.801
$ws.Range("L132").Value = 500006220
$ws.Range("M132").Value = -26609.801
$ws.Range("N132").Value = -500011280
$ws.Range("H136").Value = 20412442
$ws.Range("I136").Value = 24394088
$ws.Range("J136").Value = 6499.25
$ws.Range("K136").Value = 73182264
$ws.Range("L136").Value = 19497.75
$ws.Range("M136").Value = -73179714
$ws.Range("N136").Value = -24597.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 54856
$ws.Range("J13").Value = 54856
$ws.Range("L13").Value = 54856
$ws.Range("N13").Value = -55192
$ws.Range("H22").Value = 314.5
$ws.Range("I22").Value = 236
$ws.Range("K22").Value = 236
$ws.Range("M22").Value = -63
$ws.Range("H94").Value = 3013.5833
$ws.Range("I94").Value = 3455.625
$ws.Range("J94").Value = 2129.5
$ws.Range("K94").Value = 3455.625
$ws.Range("L94").Value = 2129.5
$ws.Range("M94").Value = -3004.625
$ws.Range("N94").Value = -3031.5
$ws.Range("H105").Value = 17085.076
$ws.Range("I105").Value = 21720.6
$ws.Range("K105").Value = 21720.6
$ws.Range("M105").Value = -19973.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 11596.556
$ws.Range("I22").Value = 16905.666
$ws.Range("K22").Value = 16905.666
$ws.Range("M22").Value = -16555.666
$ws.Range("H41").Value = 29815
$ws.Range("I41").Value = 18186.334
$ws.Range("J41").Value = 36792.2
$ws.Range("K41").Value = 18186.334
$ws.Range("L41").Value = 36792.2
$ws.Range("M41").Value = -17758.334
$ws.Range("N41").Value = -37648.2
$ws.Range("H51").Value = 49995
$ws.Range("I51").Value = 49995
$ws.Range("K51").Value = 49995
$ws.Range("M51").Value = -49259
$ws.Range("H52").Value = 149989
$ws.Range("J52").Value = 149989
$ws.Range("L52").Value = 149989
$ws.Range("N52").Value = -150577
$ws.Range("H58").Value = 3557.9583
$ws.Range("I58").Value = 3684.4736
$ws.Range("J58").Value = 3077.2
$ws.Range("K58").Value = 3684.4736
$ws.Range("L58").Value = 3077.2
$ws.Range("M58").Value = -3481.4736
$ws.Range("N58").Value = -3483.2
$ws.Range("H61").Value = 49995
$ws.Range("I61").Value = 49995
$ws.Range("K61").Value = 49995
$ws.Range("M61").Value = -49647
$ws.Range("H99").Value = 12684.1
$ws.Range("I99").Value = 17491.166
$ws.Range("K99").Value = 17491.166
$ws.Range("M99").Value = -15993.166
$ws.Range("H107").Value = 2855.2942
$ws.Range("I107").Value = 2659.2307
$ws.Range("K107").Value = 2659.2307
$ws.Range("M107").Value = -739.2307000000001
$ws.Range("H126").Value = 12684.1
$ws.Range("I126").Value = 17491.166
$ws.Range("K126").Value = 52473.49800000001
$ws.Range("M126").Value = -50003.49800000001
$ws.Range("H132").Value = 80230.234
$ws.Range("I132").Value = 93361.13
$ws.Range("K132").Value = 280083.39
$ws.Range("M132").Value = -277553.39
$ws.Range("H136").Value = 3557.9583
$ws.Range("I136").Value = 3684.4736
$ws.Range("J136").Value = 3077.2
$ws.Range("K136").Value = 11053.4208
$ws.Range("L136").Value = 9231.599999999999
$ws.Range("M136").Value = -8503.4208
$ws.Range("N136").Value = -14331.6
$ws.Range("H139").Value = 140000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 142.5
$ws.Range("I17").Value = 85
$ws.Range("K17").Value = 255
$ws.Range("M17").Value = -86
$ws.Range("H131").Value = 1903.8334
$ws.Range("J131").Value = 1896.2222
$ws.Range("L131").Value = 5688.6666
$ws.Range("N131").Value = -15768.6666
$ws.Range("H133").Value = 8889
$ws.Range("I133").Value = 7301.7144
$ws.Range("K133").Value = 21905.1432
$ws.Range("M133").Value = -16845.1432
$ws.Range("H134").Value = 7173.2964
$ws.Range("I134").Value = 1983.1578
$ws.Range("J134").Value = 19499.875
$ws.Range("K134").Value = 5949.4734
$ws.Range("L134").Value = 58499.625
$ws.Range("M134").Value = -879.4733999999999
$ws.Range("N134").Value = -68639.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 109725.266
$ws.Range("I70").Value = 137793.73
$ws.Range("K70").Value = 137793.73
$ws.Range("M70").Value = -137523.73
$ws.Range("H73").Value = 109725.266
$ws.Range("I73").Value = 137793.73
$ws.Range("K73").Value = 137793.73
$ws.Range("M73").Value = -136857.73
$ws.Range("H80").Value = 6977.7144
$ws.Range("J80").Value = 6977.7144
$ws.Range("L80").Value = 6977.7144
$ws.Range("N80").Value = -8973.714400000001
$ws.Range("H83").Value = 6977.7144
$ws.Range("J83").Value = 6977.7144
$ws.Range("L83").Value = 34888.572
$ws.Range("N83").Value = -44872.572
$ws.Range("H132").Value = 5638.6587
$ws.Range("J132").Value = 7279
$ws.Range("L132").Value = 21837
$ws.Range("N132").Value = -26897

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1769.7368
$ws.Range("J46").Value = 4480
$ws.Range("L46").Value = 4480
$ws.Range("N46").Value = -4856
$ws.Range("H82").Value = 2086.8333
$ws.Range("J82").Value = 5000
$ws.Range("L82").Value = 5000
$ws.Range("N82").Value = -5722
$ws.Range("H85").Value = 2086.8333
$ws.Range("J85").Value = 5000
$ws.Range("L85").Value = 5000
$ws.Range("N85").Value = -7496
$ws.Range("H93").Value = 3894.2
$ws.Range("I93").Value = 1836
$ws.Range("K93").Value = 1836
$ws.Range("M93").Value = -588
$ws.Range("H129").Value = 27000
$ws.Range("I129").Value = 27000
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 27000
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -22000
$ws.Range("N129").ClearContents()
$ws.Range("H132").Value = 90914270
$ws.Range("I132").Value = 5098.4165
$ws.Range("J132").Value = 200005280
$ws.Range("K132").Value = 15295.2495
$ws.Range("L132").Value = 600015840
$ws.Range("M132").Value = -12765.2495
$ws.Range("N132").Value = -600020900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6254.4897
$ws.Range("I132").Value = 6129.951
$ws.Range("K132").Value = 18389.853
$ws.Range("M132").Value = -15859.853
$ws.Range("H136").Value = 1363.1666
$ws.Range("I136").Value = 1377.579
$ws.Range("K136").Value = 4132.737
$ws.Range("M136").Value = -1582.737

